# Adds the fourth slide ("Hackerman") using the Title Only layout, with
# a title placeholder and two body textboxes (code listings), matching the
# authored OOXML exactly (text, run/paragraph structure, size & position).

$p = $ppt.ActivePresentation

# --- New slide 4, "Title Only" layout (slideLayout6.xml) ---------------
$s = $p.Slides.Add(4, 6)

# --- Title placeholder ---------------------------------------------------
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Hackerman"
$title.TextFrame.TextRange.Font.Size = 72
$title.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$title.Left = 67.8265380859375
$title.Top = 23.241655349731445
$title.Width = 425.25
$title.Height = 94.45157623291016

# --- TextBox 2: "Machine Object" code listing -----------------------------
$tb2 = $s.Shapes.AddTextbox(1, 100, 100, 200, 200)
$tb2.Name = "TextBox 2"
$tb2.TextFrame.WordWrap = -1
$tb2.TextFrame.AutoSize = 1
$tb2.Fill.Visible = 0
$tb2.TextFrame.TextRange.Text = "Machine Object`r(define (machine ip ports protocols)`r   (define open-tcp '())`r   (define open-udp '())`r   (define udp-socket (udp-open-socket) )`r   (define msg (make-bytes 100) )`r   (define (add-udp port)`r    `t(set! open-udp (cons port open-udp)))`r   (define (add-tcp port)`r        (set! open-tcp (cons port open-tcp)))`r   (define (check-uport port)`r        (if (memq port open-udp) #t #f))`r   (define (check-tport port)`r        (if (memq (string->number port) open-           `t                    tcp) #t #f))"
$tb2.TextFrame.TextRange.Paragraphs(1, 1).Font.Size = 24
$tb2.TextFrame.TextRange.Paragraphs(1, 1).ParagraphFormat.Alignment = 2
$tb2.Left = 110.32252502441406
$tb2.Top = 117.69322967529297
$tb2.Width = 356.5161437988281
$tb2.Height = 341.7047424316406

# --- TextBox 3: "Dispatch" code listing ------------------------------------
$tb3 = $s.Shapes.AddTextbox(1, 100, 100, 200, 200)
$tb3.Name = "TextBox 3"
$tb3.TextFrame.WordWrap = -1
$tb3.TextFrame.AutoSize = 1
$tb3.Fill.Visible = 0
$tb3.TextFrame.TextRange.Text = "Dispatch`r(define (dispatch message)`r    (cond((eq? (car message) 'tports) `t`t  `t`t(match-ports open-tcp))`r         ((eq? (car message) 'dtports) open-tcp)`r         ((eq? (car message) 'uports) open-udp)`r         ((eq? (car message) 'ip) ip)`r         ((eq? (car message) 'tport) (check-tport (cadr message)) )`r         ((eq? (car message) 'uport) (check-uport (cadr message)) )`r         (else error `"Bad moves, dude`")))`r  (begin (udp-bind! udp-socket `"127.0.0.1`" 0)`r         (map (lambda (x) (probe-tcp x)) (enum-ports ports))`r         (map (lambda (x) (probe-udp x)) (enum-ports ports)) dispatch))"
$tb3.TextFrame.TextRange.Paragraphs(1, 1).Font.Size = 24
$tb3.TextFrame.TextRange.Paragraphs(1, 1).ParagraphFormat.Alignment = 2
$tb3.Left = 118.45149993896484
$tb3.Top = 481.2088317871094
$tb3.Width = 340.25811767578125
$tb3.Height = 450.7593994140625
